# Add a "Save" column (H) to the s_vals worksheet, mirroring the style
# used by the existing header cells and populating the per-row save flag.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 - same header style as the rest of row 1 (e.g. G1 "sum")
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Data values for H2:H11 taken from the source diff
$saveValues = @(0, 1, 0, 1, 1, 0, 1, 0, 0, 0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
